$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 755
$ws1.Range("F3").Value = 14186
$ws1.Range("F4").Value = 14186
$ws1.Range("F5").Value = 14275
$ws1.Range("F6").Value = 1364
$ws1.Range("F7").Value = 1398
$ws1.Range("F8").Value = 5877
$ws1.Range("F9").Value = 983
$ws1.Range("F10").Value = 572
$ws1.Range("F14").Value = 1541
$ws1.Range("F15").Value = 435
$ws1.Range("F17").Value = 1196
$ws1.Range("F18").Value = 1823
$ws1.Range("F19").Value = 914
$ws1.Range("F20").Value = 35
$ws1.Range("F21").Value = 2275
$ws1.Range("F23").Value = 810
$ws1.Range("F24").Value = 3314
$ws1.Range("F27").Value = 2389
$ws1.Range("F28").Value = 590
$ws1.Range("F31").Value = 1783
$ws1.Range("F32").Value = 1068
$ws1.Range("F33").Value = 1386
$ws1.Range("F34").Value = 99
$ws1.Range("F35").Value = 148
$ws1.Range("F36").Value = 4794
$ws1.Range("F37").Value = 4837
$ws1.Range("F38").Value = 301
$ws1.Range("F39").Value = 159
$ws1.Range("F41").Value = 683
$ws1.Range("F42").Value = 3285
$ws1.Range("F43").Value = 43
$ws1.Range("F44").Value = 922
$ws1.Range("F45").Value = 334
$ws1.Range("F46").Value = 102
$ws1.Range("F47").Value = 76
$ws1.Range("F48").Value = 4419
$ws1.Range("F49").Value = 578
$ws1.Range("F50").Value = 290

# Sheet: 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F15").Value = 0
$ws2.Range("F20").Value = 13
$ws2.Range("F22").Value = 55

# Sheet: 本地生活
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 7553
$ws3.Range("F3").Value = 238
$ws3.Range("F4").Value = 764

# Sheet: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 7553
$ws4.Range("F3").Value = 755
$ws4.Range("F4").Value = 238
$ws4.Range("F5").Value = 764
$ws4.Range("F6").Value = 14186
$ws4.Range("F7").Value = 14275
$ws4.Range("F8").Value = 1364
$ws4.Range("F9").Value = 1398
$ws4.Range("F10").Value = 5877
$ws4.Range("F11").Value = 983
$ws4.Range("F15").Value = 1541
$ws4.Range("F16").Value = 435
$ws4.Range("F17").Value = 1196
$ws4.Range("F18").Value = 1823
$ws4.Range("F19").Value = 914
$ws4.Range("F20").Value = 35
$ws4.Range("F21").Value = 3314
$ws4.Range("F23").Value = 2389
$ws4.Range("F24").Value = 590
$ws4.Range("F27").Value = 1783
$ws4.Range("F29").Value = 20
$ws4.Range("F31").Value = 1068
$ws4.Range("F32").Value = 1386
$ws4.Range("F33").Value = 99
$ws4.Range("F34").Value = 4794
$ws4.Range("F35").Value = 4837
$ws4.Range("F36").Value = 301
$ws4.Range("F37").Value = 159
$ws4.Range("F39").Value = 683
$ws4.Range("F40").Value = 3285
$ws4.Range("F41").Value = 922
$ws4.Range("F42").Value = 334
$ws4.Range("F43").Value = 102
$ws4.Range("F45").Value = 76
$ws4.Range("F46").Value = 4419
$ws4.Range("F47").Value = 578
$ws4.Range("F48").Value = 290
